$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows before the old row 16 (pushes the "storage table" block
# from rows 16-22 down to rows 23-29, formulas/refs shift automatically).
$ws.Range("A15:A21").EntireRow.Insert()

# --- New content, written in the same order the original author typed it
# (this controls the order new entries land in xl/sharedStrings.xml) ---

$ws.Range("A31").Value = "Matrice G"

$ws.Range("B16").Value = "Matrix G"
$ws.Range("B17").Value = "Matrix B"
$ws.Range("B18").Value = "Vector s"

$ws.Range("A34").Value = "Matrix B and vector s"

$ws.Range("B19").Value = "Matrix A"

$ws.Range("A38").Value = "Benoit Foret"
$ws.Range("B38").Value = "4TB"

$ws.Range("C15").Value = "Direct object size"

# --- Remaining values / labels in the freshly inserted rows (15-21) ---

$ws.Range("A16").Value = "per iter and act"
$ws.Range("C16").Value = 2.5

$ws.Range("A17").Value = "per iteration"
$ws.Range("C17").Value = 5

$ws.Range("A18").Value = "per iter and act"
$ws.Range("C18").Value = 0.109

$ws.Range("A19").Value = "per iteration"
$ws.Range("C19").Value = 2.5

# --- New content below the (now shifted) TOTAL row (row 29) ---

$ws.Range("A32").Value = "per iter and act"
$ws.Range("C32").Formula = "=B28*C16"
$ws.Range("D32").Formula = "=C32/1000"
$ws.Range("E32").Formula = "=D32/1000"
$ws.Range("E32").NumberFormat = "0.00"

$ws.Range("A35").Value = "per iteration"
$ws.Range("B35").Formula = "=B27"
$ws.Range("C35").Formula = "=B35*(C17+C19)"
$ws.Range("D35").Formula = "=C35/1000"
$ws.Range("E35").Formula = "=D35/1000"
$ws.Range("E35").NumberFormat = "0.00"

$ws.Range("A36").Value = "per iter and act"
$ws.Range("B36").Formula = "=B28"
$ws.Range("C36").Formula = "=B36*C18"
$ws.Range("D36").Formula = "=C36/1000"
$ws.Range("E36").Formula = "=D36/1000"
$ws.Range("E36").NumberFormat = "0.00"

# --- View / selection bookkeeping to match the authored workbook ---

$ws.Range("H30").NumberFormat = "0.00"

$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Range("C16").Select()
